$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the three input cells (price values) from 1 to their new amounts
$ws.Range("I3").Value = 500
$ws.Range("I8").Value = 500
$ws.Range("I13").Value = 100

# Update the active cell selection to match the final state of the file
$ws.Activate()
$ws.Range("I14").Select()
